# "Add files via upload" — the uploaded workbook gained a new data row/cell
# (C2 = "d") and the previously-entered cells (B1, C1, A2, B2) picked up a
# freshly-minted style entry (a clone of the original font/format) while A1
# was left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value: C2 = "d" (becomes a new shared-string entry).
$ws.Cells.Item(2, 3).Value = "d"

# Touch the formatting of every other populated cell except A1 so Excel
# clones a new style/font entry for them (A1 keeps its original style).
$xlThemeColorDark1 = 1
$ws.Range("B1:C1").Font.ThemeColor = $xlThemeColorDark1
$ws.Range("A2:C2").Font.ThemeColor = $xlThemeColorDark1
